# This edit re-syncs the "Artfynd" sheet with an updated export where several
# observation records (rows) were re-sequenced. For each affected pair of rows
# the two records fully swap their field values (Id, Ost/Nord coordinates,
# activity/comment/substrate texts, species, etc.) while the row position
# itself stays fixed. Cells that become blank are cleared; cells that gain a
# value are set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (was record 131112844, becomes 131112845)
$ws.Range("A4").Value = 131112845
$ws.Range("M4").Value = 'färska spår'
$ws.Range("Q4").Value = 510982
$ws.Range("R4").Value = 7037550
$ws.Range("AC4").Value = 'Ringhack, färska, på gran.'
$ws.Range("AM4").Value = 'Trädstam på levande träd'
$ws.Range("AO4").Value = 'Stem on living tree # Picea abies'

# Row 5 (was record 131112845, becomes 131112844)
$ws.Range("A5").Value = 131112844
$ws.Range("M5").Value = 'äldre spår'
$ws.Range("Q5").Value = 510994
$ws.Range("R5").Value = 7037511
$ws.Range("AC5").Value = 'Ringhack, äldre, på gran.'
$ws.Range("AM5").ClearContents()
$ws.Range("AO5").Value = 'Picea abies'

# Row 11 (was record 131112847, becomes 131112860)
$ws.Range("A11").Value = 131112860
$ws.Range("B11").Value = 79243
$ws.Range("E11").Value = 6425
$ws.Range("F11").Value = 'Garnlav'
$ws.Range("G11").Value = 'Alectoria sarmentosa'
$ws.Range("H11").Value = '(Ach.) Ach.'
$ws.Range("J11").ClearContents()
$ws.Range("L11").ClearContents()
$ws.Range("M11").ClearContents()
$ws.Range("Q11").Value = 511032
$ws.Range("R11").Value = 7037541
$ws.Range("AC11").ClearContents()
$ws.Range("AF11").ClearContents()

# Row 12 (was record 131112860, becomes 131112847)
$ws.Range("A12").Value = 131112847
$ws.Range("B12").Value = 57884
$ws.Range("E12").Value = 100109
$ws.Range("F12").Value = 'Tretåig hackspett'
$ws.Range("G12").Value = 'Picoides tridactylus'
$ws.Range("H12").Value = '(Linnaeus, 1758)'
$ws.Range("J12").ClearContents()
$ws.Range("L12").ClearContents()
$ws.Range("M12").Value = 'äldre spår'
$ws.Range("Q12").Value = 511011
$ws.Range("R12").Value = 7037468
$ws.Range("AC12").Value = 'Ringhack, äldre, på gran.'
$ws.Range("AF12").ClearContents()

# Row 17 (was record 131112855, becomes 131112852)
$ws.Range("A17").Value = 131112852
$ws.Range("B17").Value = 57881
$ws.Range("E17").Value = 100049
$ws.Range("F17").Value = 'Spillkråka'
$ws.Range("G17").Value = 'Dryocopus martius'
$ws.Range("H17").Value = '(Linnaeus, 1758)'
$ws.Range("J17").ClearContents()
$ws.Range("L17").ClearContents()
$ws.Range("M17").Value = 'äldre spår'
$ws.Range("Q17").Value = 511154
$ws.Range("R17").Value = 7037711
$ws.Range("AC17").Value = 'Rejäla äldre hackspår i stambasen av en gran.'
$ws.Range("AF17").ClearContents()
$ws.Range("AJ17").Value = 'gran'
$ws.Range("AK17").Value = 'Picea abies'
$ws.Range("AM17").Value = 'Stående död trädstam/högstubbe'
$ws.Range("AO17").Value = 'Standing dead tree/snags # Picea abies'

# Row 18 (was record 131112852, becomes 131112855)
$ws.Range("A18").Value = 131112855
$ws.Range("B18").Value = 80348
$ws.Range("E18").Value = 6458
$ws.Range("F18").Value = 'Lunglav'
$ws.Range("G18").Value = 'Lobaria pulmonaria'
$ws.Range("H18").Value = '(L.) Hoffm.'
$ws.Range("J18").ClearContents()
$ws.Range("L18").ClearContents()
$ws.Range("M18").ClearContents()
$ws.Range("Q18").Value = 511045
$ws.Range("R18").Value = 7037478
$ws.Range("AC18").ClearContents()
$ws.Range("AF18").ClearContents()
$ws.Range("AJ18").Value = 'sälg'
$ws.Range("AK18").Value = 'Salix caprea'
$ws.Range("AM18").Value = 'Bark på levande träd'
$ws.Range("AO18").Value = 'Bark on living woody plant # Salix caprea'

# Row 19 (was record 131112856, becomes 131112848)
$ws.Range("A19").Value = 131112848
$ws.Range("B19").Value = 57884
$ws.Range("E19").Value = 100109
$ws.Range("F19").Value = 'Tretåig hackspett'
$ws.Range("G19").Value = 'Picoides tridactylus'
$ws.Range("H19").Value = '(Linnaeus, 1758)'
$ws.Range("J19").ClearContents()
$ws.Range("K19").ClearContents()
$ws.Range("L19").ClearContents()
$ws.Range("M19").Value = 'äldre spår'
$ws.Range("Q19").Value = 510963
$ws.Range("R19").Value = 7037437
$ws.Range("AC19").Value = 'Ringhack, äldre, på gran.'
$ws.Range("AF19").ClearContents()

# Row 20 (was record 131112848, becomes 131112856)
$ws.Range("A20").Value = 131112856
$ws.Range("B20").Value = 91804
$ws.Range("E20").Value = 1108
$ws.Range("F20").Value = 'Harticka'
$ws.Range("G20").Value = 'Pelloporus leporinus'
$ws.Range("H20").Value = '(Fr.) Krieglst.'
$ws.Range("J20").ClearContents()
$ws.Range("K20").Value = 'teleomorf'
$ws.Range("L20").ClearContents()
$ws.Range("M20").ClearContents()
$ws.Range("Q20").Value = 511011
$ws.Range("R20").Value = 7037561
$ws.Range("AC20").ClearContents()
$ws.Range("AF20").ClearContents()

# Row 21 (was record 131112842, becomes 131112843)
$ws.Range("A21").Value = 131112843
$ws.Range("Q21").Value = 510984
$ws.Range("R21").Value = 7037439

# Row 22 (was record 131112843, becomes 131112842)
$ws.Range("A22").Value = 131112842
$ws.Range("Q22").Value = 511092
$ws.Range("R22").Value = 7037598
